$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2, shifting existing data rows (2-13) down to (3-14)
$ws.Rows.Item(2).Insert()

# The insert picks up formatting from the row above (the bold header row);
# clear it back to the plain/default look used by the rest of the data rows.
$ws.Rows.Item(2).ClearFormats()

# Re-apply the date number format used by the other "Fecha" cells in column D
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new row 2 with the new record's data
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44812
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100112036
$ws.Range("G2").Value = "Caigua"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 20000
$ws.Range("N2").Value = "$/caja 15 kilos"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 1333
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = "Hortaliza"
